# Rename headers and fix municipality/state name capitalization
# (lowercase connector words "de/del/el/los/las/y" -> capitalized,
# and normalize "GUANAJUATO" -> "Guanajuato"), per data-cleaning fix.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'mx_state'
$ws.Range("B1").Value = 'mx_municipality'
$ws.Range("C1").Value = 'n_matriculas'
$ws.Range("D1").Value = 'pct_matriculas'
$ws.Range("B7").Value = 'Amatenango De La Frontera'
$ws.Range("B9").Value = 'Bejucal De Ocampo'
$ws.Range("B14").Value = 'Comitán De Domínguez'
$ws.Range("B27").Value = 'San Cristóbal De Las Casas'
$ws.Range("A42").Value = 'Ciudad De México'
$ws.Range("B44").Value = 'Cuajimalpa De Morelos'
$ws.Range("A54").Value = 'Estado De México'
$ws.Range("B54").Value = 'Almoloya De Juárez'
$ws.Range("B61").Value = 'Ecatepec De Morelos'
$ws.Range("B65").Value = 'Naucalpan De Juárez'
$ws.Range("B67").Value = 'San Felipe Del Progreso'
$ws.Range("B71").Value = 'Tlalnepantla De Baz'
$ws.Range("A78").Value = 'Guanajuato'
$ws.Range("B78").Value = 'Apaseo El Alto'
$ws.Range("B91").Value = 'Acapulco De Juárez'
$ws.Range("B93").Value = 'Alcozauca De Guerero'
$ws.Range("B96").Value = 'Atoyac De Álvarez'
$ws.Range("B98").Value = 'Chilpancingo De Los Bravo'
$ws.Range("B99").Value = 'Coahuayutla De José María Izazaga'
$ws.Range("B103").Value = 'Coyuca De Benítez'
$ws.Range("B107").Value = 'Iguala De La Independencia'
$ws.Range("B109").Value = 'Zihuatanejo De Azueta'
$ws.Range("B118").Value = 'Técpan De Galeana'
$ws.Range("B120").Value = 'Tixtla De Guerero'
$ws.Range("B121").Value = 'Tlapa De Comonfort'
$ws.Range("B137").Value = 'Mixquiahuala De Juárez'
$ws.Range("B139").Value = 'Pachuca De Soto'
$ws.Range("B148").Value = 'Autlán De Navarro'
$ws.Range("B151").Value = 'Encarnación De Díaz'
$ws.Range("B154").Value = 'Ixtlahuacán Del Río'
$ws.Range("B155").Value = 'Jilotlán De Los Dolores'
$ws.Range("B158").Value = 'Lagos De Moreno'
$ws.Range("B164").Value = 'San Juanito De Escobedo'
$ws.Range("B165").Value = 'San Miguel El Alto'
$ws.Range("B166").Value = 'Teocuitatlán De Corona'
$ws.Range("B167").Value = 'Tizapán El Alto'
$ws.Range("B171").Value = 'Valle De Juárez'
$ws.Range("B197").Value = 'Tiquicheo De Nicolás Romero'
$ws.Range("B209").Value = 'Jonacatepec De Leandro Valle'
$ws.Range("B213").Value = 'Tetela Del Volcán'
$ws.Range("B217").Value = 'Zacualpan De Amilpas'
$ws.Range("B224").Value = 'Acatlán De Pérez Figueroa'
$ws.Range("B228").Value = 'Fresnillo De Trujano'
$ws.Range("B229").Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range("B230").Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range("B232").Value = 'Nejapa De Madero'
$ws.Range("B233").Value = 'Oaxaca De Juárez'
$ws.Range("B234").Value = 'Ocotlán De Morelos'
$ws.Range("B235").Value = 'Putla Villa De Guerero'
$ws.Range("B239").Value = 'San Antonino El Alto'
$ws.Range("B260").Value = 'San Miguel El Grande'
$ws.Range("B264").Value = 'San Pedro Y San Pablo Teposcolula'
$ws.Range("B267").Value = 'Santa Ana Del Valle'
$ws.Range("B278").Value = 'Santo Domingo De Morelos'
$ws.Range("B281").Value = 'Tlalixtac De Cabrera'
$ws.Range("B282").Value = 'Totontepec Villa De Morelos'
$ws.Range("B283").Value = 'Villa De Zaachila'
$ws.Range("B284").Value = 'Villa Sola De Vega'
$ws.Range("B285").Value = 'Zimatlán De Álvarez'
$ws.Range("B302").Value = 'Chila De La Sal'
$ws.Range("B307").Value = 'Cuapiaxtla De Madero'
$ws.Range("B310").Value = 'Cuayuca De Andrade'
$ws.Range("B320").Value = 'Huehuetlán El Chico'
$ws.Range("B321").Value = 'Huehuetlán El Grande'
$ws.Range("B323").Value = 'Ixcamilpa De Guerero'
$ws.Range("B326").Value = 'Izúcar De Matamoros'
$ws.Range("B331").Value = 'Los Reyes De Juárez'
$ws.Range("B351").Value = 'San Nicolás De Los Ranchos'
$ws.Range("B355").Value = 'San Salvador El Seco'
$ws.Range("B358").Value = 'Tecali De Herrera'
$ws.Range("B366").Value = 'Tepexi De Rodríguez'
$ws.Range("B368").Value = 'Tepeyahualco De Cuauhtémoc'
$ws.Range("B371").Value = 'Tlacotepec De Benito Juárez'
$ws.Range("B380").Value = 'Xayacatlán De Bravo'
$ws.Range("B391").Value = 'Pinal De Amoles'
$ws.Range("B393").Value = 'San Juan Del Río'
$ws.Range("B398").Value = 'Ciudad Del Maíz'
$ws.Range("B402").Value = 'Santa María Del Río'
$ws.Range("B429").Value = 'Contla De Juan Cuamatzi'
$ws.Range("B432").Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range("B434").Value = 'Mazatecochco De José María Morelos'
$ws.Range("B435").Value = 'Nanacamilpa De Mariano Arista'
$ws.Range("B437").Value = 'Papalotla De Xicohténcatl'
$ws.Range("B441").Value = 'San Pablo Del Monte'
$ws.Range("B445").Value = 'Tepetitla De Lardizábal'
$ws.Range("B447").Value = 'Tetla De La Solidaridad'
$ws.Range("B466").Value = 'Ignacio De La Llave'
$ws.Range("B468").Value = 'Ixhuatlán De Madero'
$ws.Range("B475").Value = 'Martínez De La Torre'

# Fix two floating point rounding values (re-derived during recompute)
$ws.Range("D28").Value = 0.009925558312655089
$ws.Range("D286").Value = 0.09243176178660048

# Remove trailing footer/metadata rows (504-508): sample size, source,
# author and date notes that are not part of the tabular data.
$ws.Range("A504:D508").EntireRow.Delete()

